# ---------------------------------------------------------------------------
# "Add files via upload" — refresh the measurement data (columns B:C) on
# Sheet1, let Excel recompute the dependent T-column running sums and the
# chart's cached values, reposition the embedded chart, and update the
# active selection / scroll position.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$data = New-Object 'object[,]' 40,2
$data[0,0] = 940625
$data[0,1] = 990320
$data[1,0] = 1895511
$data[1,1] = 2029508
$data[2,0] = 2828932
$data[2,1] = 3187054
$data[3,0] = 3830703
$data[3,1] = 4549176
$data[4,0] = 5107226
$data[4,1] = 6084360
$data[5,0] = 6238901
$data[5,1] = 7514827
$data[6,0] = 7239661
$data[6,1] = 9145456
$data[7,0] = 8942810
$data[7,1] = 11323266
$data[8,0] = 9431490
$data[8,1] = 12671432
$data[9,0] = 11691806
$data[9,1] = 15038477
$data[10,0] = 12527370
$data[10,1] = 16588047
$data[11,0] = 14479240
$data[11,1] = 19027607
$data[12,0] = 15725580
$data[12,1] = 21749021
$data[13,0] = 16177867
$data[13,1] = 23101046
$data[14,0] = 18265559
$data[14,1] = 25747753
$data[15,0] = 19623979
$data[15,1] = 27889815
$data[16,0] = 20682123
$data[16,1] = 30590036
$data[17,0] = 21362086
$data[17,1] = 32641962
$data[18,0] = 22667913
$data[18,1] = 35741688
$data[19,0] = 27200803
$data[19,1] = 39934480
$data[20,0] = 26206265
$data[20,1] = 41220365
$data[21,0] = 29652915
$data[21,1] = 45258258
$data[22,0] = 31237419
$data[22,1] = 48650006
$data[23,0] = 31854566
$data[23,1] = 50620629
$data[24,0] = 34537990
$data[24,1] = 54732508
$data[25,0] = 35798192
$data[25,1] = 58215365
$data[26,0] = 39382465
$data[26,1] = 62603484
$data[27,0] = 38033011
$data[27,1] = 64120789
$data[28,0] = 42054907
$data[28,1] = 68665467
$data[29,0] = 42797005
$data[29,1] = 72898932
$data[30,0] = 44011903
$data[30,1] = 76003025
$data[31,0] = 46858856
$data[31,1] = 80160440
$data[32,0] = 48937093
$data[32,1] = 85026313
$data[33,0] = 50428823
$data[33,1] = 88389041
$data[34,0] = 53898727
$data[34,1] = 93688611
$data[35,0] = 55139263
$data[35,1] = 97228870
$data[36,0] = 55829868
$data[36,1] = 100748836
$data[37,0] = 57625595
$data[37,1] = 105644874
$data[38,0] = 61549820
$data[38,1] = 112217880
$data[39,0] = 64328131
$data[39,1] = 115697648

# Columns B (idx 0) and C (idx 1) for rows 2..41 hold the raw measurement
# sums; T3:T42 = SUM(B{r-1}:C{r-1}) and the chart series cache recompute
# automatically off of these once Excel recalculates.
$ws.Range("B2:C41").Value = $data

# Reposition the embedded chart (moved up/right and shrunk slightly).
$co = $ws.ChartObjects().Item(1)
$co.Left = 725.9999212598425
$co.Top = 120
$co.Width = 558.3486328124999
$co.Height = 267.37496062992125

# Update the view: drop the frozen "topLeftCell" scroll anchor and move the
# active selection from Q11 to X13.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("X13").Select()
